$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-9 down to 9-10
$ws.Rows.Item(8).Insert()

# Set values for the new row 8
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 44511
$ws.Cells.Item(8, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100104
$ws.Cells.Item(8, 8).Value = "Frutos de pepita"
$ws.Cells.Item(8, 9).Value = 100104005
$ws.Cells.Item(8, 10).Value = "Pera asiática"
$ws.Cells.Item(8, 11).Value = "Hosui"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 15
$ws.Cells.Item(8, 14).Value = 22000
$ws.Cells.Item(8, 15).Value = 22000
$ws.Cells.Item(8, 16).Value = 22000
$ws.Cells.Item(8, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 1467
$ws.Cells.Item(8, 20).Value = 15
